$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)

# 1. Grow the text box to make room for the two new bullet lines.
#    PowerPoint's Shape.Height/Width are expressed in points (1 pt = 12700 EMU),
#    so convert the target EMU value (2031325) to points.
$shp.Height = 2031325 / 12700

$tr = $shp.TextFrame.TextRange

# 2. Split the trailing run of paragraph 5 ("...顯示文件變更的記錄的詳細信息")
#    into two runs ("...詳細" | "信息") by re-assigning the last two
#    characters to themselves - this forces PowerPoint to break the run
#    without altering any visible formatting.
$para5 = $tr.Paragraphs(5, 1)
$para5Text = $para5.Text
$tail = $para5.Characters($para5Text.Length - 1, 2)
$tail.Text = "信息"

# 3. Append the two new bullet paragraphs after paragraph 5. Using `r
#    (carriage return) creates brand-new paragraphs that inherit the
#    same pPr (marL/indent/buAutoNum) as paragraph 5.
$tr = $shp.TextFrame.TextRange
$para5 = $tr.Paragraphs(5, 1)
$newText = $para5.InsertAfter("`rgit stash save stash信息 --include-untrack --- 将未被追踪的文件也存储起来`rgit stash branch 分支名     ---  将存储的工作区运用到新的branch")

# 4. Re-split the two freshly inserted paragraphs into the same run
#    boundaries the authored deck uses (each boundary is re-asserted via
#    a same-text write, which breaks the run without touching formatting).
function Split-Runs($paragraph, $segments) {
    $pos = 1
    foreach ($seg in $segments) {
        $len = $seg.Length
        if ($len -gt 0) {
            $chunk = $paragraph.Characters($pos, $len)
            $chunk.Text = $seg
        }
        $pos += $len
    }
}

$tr = $shp.TextFrame.TextRange
$para6 = $tr.Paragraphs(6, 1)
Split-Runs $para6 @("git", " stash save ", "stash", "信息", " ", "--", "include-untrack --- ", "将未被追踪的文件也存储起来")

$tr = $shp.TextFrame.TextRange
$para7 = $tr.Paragraphs(7, 1)
Split-Runs $para7 @("git", " stash branch ", "分支名", "     ---  ", "将存储的工作区运用到新的", "branch")

Write-Host "Paragraph count: $($shp.TextFrame.TextRange.Paragraphs().Count)"
Write-Host $shp.TextFrame.TextRange.Text
